# Populate the two empty cells in row "5" of the "ХРОНОЛОГИЯ НА СЪБИТИЯТА"
# (event chronology) table with the company name and the date/time of the
# correspondence, matching the previously-filled rows above it.

$d = $word.ActiveDocument

# Locate the chronology table: it's the 3-column table whose second
# header row starts with "N" (the "N / Компания / Дата и час ..." row).
# Fall back to the last table in the document if that search fails.
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    if ($tbl.Columns.Count -eq 3 -and $tbl.Rows.Count -ge 2) {
        $headerText = $tbl.Cell(2, 1).Range.Text
        if ($headerText -like "N*") {
            $targetTable = $tbl
        }
    }
}
if ($targetTable -eq $null) {
    $targetTable = $d.Tables.Item($d.Tables.Count)
}

# Find the data row whose first (numbering) column holds "5" (cell text
# always carries a trailing cell-mark, hence the wildcard match) and
# whose company/date columns are still empty.
$targetRow = -1
for ($r = 1; $r -le $targetTable.Rows.Count; $r++) {
    $numText = $targetTable.Cell($r, 1).Range.Text
    $companyText = $targetTable.Cell($r, 2).Range.Text
    if ($numText -like "5*" -and $companyText.Length -le 2) {
        $targetRow = $r
    }
}

# Only act if a matching (still-empty) row was actually found, so the
# script is a safe no-op rather than clobbering an unrelated row if the
# document doesn't look the way we expect.
if ($targetRow -ge 1) {
    # Fill the "Компания" (company) cell.
    $companyCell = $targetTable.Cell($targetRow, 2)
    $companyCell.Range.Text = "CHINA INC."
    $targetTable.Cell($targetRow, 2).Range.Font.Name = "Calibri"

    # Fill the "Дата и час на кореспонденция" (date/time) cell.
    $dateCell = $targetTable.Cell($targetRow, 3)
    $dateCell.Range.Text = "20.10.2021 / 18:52"
    $targetTable.Cell($targetRow, 3).Range.Font.Name = "Calibri"
}
